# The author opened the "Heating" sheet, zoomed in to take a closer look
# at the small grid, and clicked into the empty cell E2 (just to the right
# of the last populated cell in row 2) - touching its formatting without
# ever entering a value, which is enough for Excel to materialize the
# cell with its own style entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heating")
$ws.Activate()
$excel.ActiveWindow.Zoom = 169

$cell = $ws.Range("E2")
$cell.Select()
$cell.Font.ThemeColor = 1
